$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.723.89"
$ws.Range("E2").Value = "  +0.92%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.151.55"
$ws.Range("E3").Value = "  +1.89%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.10"
$ws.Range("E5").Value = "  +2.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.97"
$ws.Range("E6").Value = "  +4.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.150.80"
$ws.Range("E8").Value = "  +1.93%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.527"
$ws.Range("E9").Value = "  +4.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.162"
$ws.Range("E10").Value = "  +5.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.17"
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.503"
$ws.Range("E12").Value = "  +6.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000255"
$ws.Range("E13").Value = "  +11.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.75"
$ws.Range("E14").Value = "  +7.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.663.68"
$ws.Range("E15").Value = "  +1.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.821.60"
$ws.Range("E16").Value = "  +0.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.23"
$ws.Range("E17").Value = "  +6.83%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.157.80"
$ws.Range("E18").Value = "  +2.22%  "
$ws.Range("E19").Value = "  +0.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "511.71"
$ws.Range("E20").Value = "  +5.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.96"
$ws.Range("E21").Value = "  +7.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.734"
$ws.Range("E22").Value = "  +8.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.16"
$ws.Range("E23").Value = "  +8.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.83"
$ws.Range("E24").Value = "  +3.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.16"
$ws.Range("E25").Value = "  +4.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("E27").Value = "  +4.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.70"
$ws.Range("E28").Value = "  +8.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.19"
$ws.Range("E29").Value = "  +5.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "28.01"
$ws.Range("E30").Value = "  +6.38%  "
$ws.Range("E31").Value = "  +0.19%  "
$ws.Range("E32").Value = "  +3.58%  "
$ws.Range("E33").Value = "  +6.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.09"
$ws.Range("E34").Value = "  +9.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.56"
$ws.Range("E35").Value = "  +5.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.65"
$ws.Range("E36").Value = "  -0.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "482.45"
$ws.Range("E37").Value = "  +5.99%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0864"
$ws.Range("E38").Value = "  +5.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0422"
$ws.Range("E39").Value = "  +3.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.01"
$ws.Range("E40").Value = "  +2.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.113.04"
$ws.Range("E41").Value = "  +4.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.64"
$ws.Range("E42").Value = "  +4.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.120"
$ws.Range("E43").Value = "  +4.89%  "
$ws.Range("E44").Value = "  +12.79%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.46"
$ws.Range("E45").Value = "  +15.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.18"
$ws.Range("E46").Value = "  +4.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₃0579"
$ws.Range("E47").Value = "  +12.28%  "
$ws.Range("E49").Value = "  +3.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.29"
$ws.Range("E50").Value = "  +10.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "118.79"
$ws.Range("E51").Value = "  -1.40%  "
